# Insert a new weekly record at row 196 (Haba, Vega Central Mapocho de Santiago)
# and push the existing rows 196:265 down to 197:266.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("196:196").Insert()

$ws.Range("A196").Value = 9
$ws.Range("B196").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C196").Value = "Metropolitana"
$ws.Range("D196").Value = 44809
$ws.Range("E196").Value = 13
$ws.Range("F196").Value = 100112026
$ws.Range("G196").Value = "Haba"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 30
$ws.Range("K196").Value = 12000
$ws.Range("L196").Value = 12000
$ws.Range("M196").Value = 12000
$ws.Range("N196").Value = "`$/saco 25 kilos"
$ws.Range("O196").Value = "Provincia de Limarí"
$ws.Range("P196").Value = 480
$ws.Range("Q196").Value = 25
$ws.Range("R196").Value = "Hortaliza"
